# Rerun of models including new transit access feature.
# For each of the 16 result sheets: update coefficient values in B2:B8,
# fold the old row 11 ("Age") out by shifting names (row9->Commute_Trip,
# row10->Age) with refreshed coefficients, delete the now-duplicate row 11,
# and rename the sheet tab to its new summ<id> label.
$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 2).Value = -5.298537822659711
$ws.Cells.Item(3, 2).Value = -180.4427785809171
$ws.Cells.Item(4, 2).Value = 374.0211309922239
$ws.Cells.Item(5, 2).Value = -27.53384250283693
$ws.Cells.Item(6, 2).Value = 37.86149665076027
$ws.Cells.Item(7, 2).Value = -80.05892926866409
$ws.Cells.Item(8, 2).Value = -24.71508348881491
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 172.2444700964362
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 242.9630287414954
$ws.Rows.Item(11).Delete()
$ws.Name = "summ00639888"

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 2).Value = 7.389948162521339
$ws.Cells.Item(3, 2).Value = 25.48155765306683
$ws.Cells.Item(4, 2).Value = 146.3756012739964
$ws.Cells.Item(5, 2).Value = -4.482362206340795
$ws.Cells.Item(6, 2).Value = -259.5155030059668
$ws.Cells.Item(7, 2).Value = 9.719550362015177
$ws.Cells.Item(8, 2).Value = 4.070520661422876
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 179.0203222115342
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 17.12267731860595
$ws.Rows.Item(11).Delete()
$ws.Name = "summ00731130"

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 2).Value = 10.17572974187632
$ws.Cells.Item(3, 2).Value = 36.23806430146101
$ws.Cells.Item(4, 2).Value = 69.49976241772049
$ws.Cells.Item(5, 2).Value = 4.480140774710051
$ws.Cells.Item(6, 2).Value = -297.7475645722504
$ws.Cells.Item(7, 2).Value = -13.61939481573252
$ws.Cells.Item(8, 2).Value = 30.91020423538342
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 193.8103932864326
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 64.26119248926719
$ws.Rows.Item(11).Delete()
$ws.Name = "summ00817641"

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 2).Value = 16.16887171553983
$ws.Cells.Item(3, 2).Value = -121.1536571400732
$ws.Cells.Item(4, 2).Value = 317.9829585905956
$ws.Cells.Item(5, 2).Value = 34.96804244720039
$ws.Cells.Item(6, 2).Value = -203.9498172553326
$ws.Cells.Item(7, 2).Value = 16.07216764161825
$ws.Cells.Item(8, 2).Value = -26.28834451185726
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = -118.321833361318
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 78.24122401553024
$ws.Rows.Item(11).Delete()
$ws.Name = "summ00906655"

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 2).Value = 11.03701113706198
$ws.Cells.Item(3, 2).Value = -449.8020344105576
$ws.Cells.Item(4, 2).Value = 873.0384029006418
$ws.Cells.Item(5, 2).Value = 86.528901210296
$ws.Cells.Item(6, 2).Value = 121.8661394935399
$ws.Cells.Item(7, 2).Value = -15.64315752225758
$ws.Cells.Item(8, 2).Value = -80.65887893473848
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 67.78459816341945
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 23.71971537669236
$ws.Rows.Item(11).Delete()
$ws.Name = "summ00995476"

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 2).Value = 30.83391277188122
$ws.Cells.Item(3, 2).Value = 136.1774803469143
$ws.Cells.Item(4, 2).Value = -446.2749843690392
$ws.Cells.Item(5, 2).Value = -26.69503512618202
$ws.Cells.Item(6, 2).Value = 41.50199934591086
$ws.Cells.Item(7, 2).Value = 2.86543450187645
$ws.Cells.Item(8, 2).Value = 33.43151748694991
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 167.9318903315186
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 54.91651263761474
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01085480"

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 2).Value = 8.012236258043515
$ws.Cells.Item(3, 2).Value = -129.2237298263722
$ws.Cells.Item(4, 2).Value = 149.9711321265706
$ws.Cells.Item(5, 2).Value = 16.20844965313552
$ws.Cells.Item(6, 2).Value = 50.4316350233803
$ws.Cells.Item(7, 2).Value = 15.70593693436912
$ws.Cells.Item(8, 2).Value = -44.72759762566841
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 143.1331846454386
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 7.003806097701528
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01171988"

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 2).Value = 69.63989925904448
$ws.Cells.Item(3, 2).Value = -967.5419298869188
$ws.Cells.Item(4, 2).Value = -65.60257507042797
$ws.Cells.Item(5, 2).Value = 20.63286432866368
$ws.Cells.Item(6, 2).Value = -385.9684342071432
$ws.Cells.Item(7, 2).Value = -30.10842127571243
$ws.Cells.Item(8, 2).Value = -113.3377467219112
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 443.2934271509471
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 143.2453888361025
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01258933"

# --- Sheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 2).Value = 34.34939667794471
$ws.Cells.Item(3, 2).Value = -129.2174759795407
$ws.Cells.Item(4, 2).Value = 822.6646552828392
$ws.Cells.Item(5, 2).Value = 35.4158661000134
$ws.Cells.Item(6, 2).Value = -70.82733342514075
$ws.Cells.Item(7, 2).Value = -36.05223849243501
$ws.Cells.Item(8, 2).Value = 13.3705308793196
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 68.23003960072235
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 93.00054464481016
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01347446"

# --- Sheet 10 ---
$ws = $wb.Worksheets.Item(10)
$ws.Cells.Item(2, 2).Value = 2.21480300446407
$ws.Cells.Item(3, 2).Value = -260.4396015771994
$ws.Cells.Item(4, 2).Value = -287.7919368932203
$ws.Cells.Item(5, 2).Value = 122.8622393907799
$ws.Cells.Item(6, 2).Value = 265.7571327014712
$ws.Cells.Item(7, 2).Value = 51.50048429078774
$ws.Cells.Item(8, 2).Value = -173.1302766433438
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 133.9341232272808
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -71.51466965092541
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01435188"

# --- Sheet 11 ---
$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(2, 2).Value = 1373.459759096696
$ws.Cells.Item(3, 2).Value = 2008.179238864212
$ws.Cells.Item(4, 2).Value = 1199.533845059332
$ws.Cells.Item(5, 2).Value = 174.3708352629996
$ws.Cells.Item(6, 2).Value = 1831.48328103683
$ws.Cells.Item(7, 2).Value = -77.5662339146679
$ws.Cells.Item(8, 2).Value = -290.9334400259021
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = -50.06463841053269
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 0.5461653057452622
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01522696"

# --- Sheet 12 ---
$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(2, 2).Value = -17.26168817858459
$ws.Cells.Item(3, 2).Value = 208.6911476594801
$ws.Cells.Item(4, 2).Value = 5.725623322621317
$ws.Cells.Item(5, 2).Value = -25.33213956954594
$ws.Cells.Item(6, 2).Value = -83.77502409269675
$ws.Cells.Item(7, 2).Value = 20.97905095564724
$ws.Cells.Item(8, 2).Value = -4.487757121015534
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 174.0443843999154
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 2.786542208754469
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01613202"

# --- Sheet 13 ---
$ws = $wb.Worksheets.Item(13)
$ws.Cells.Item(2, 2).Value = -96.05035590700822
$ws.Cells.Item(3, 2).Value = 1147.057863539249
$ws.Cells.Item(4, 2).Value = 455.8684014393575
$ws.Cells.Item(5, 2).Value = 40.97162001355986
$ws.Cells.Item(6, 2).Value = 766.8283408516552
$ws.Cells.Item(7, 2).Value = 4.482716757465397
$ws.Cells.Item(8, 2).Value = -294.3922753774186
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = -88.21547334063132
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 43.70377964179605
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01700711"

# --- Sheet 14 ---
$ws = $wb.Worksheets.Item(14)
$ws.Cells.Item(2, 2).Value = 12.7023284101402
$ws.Cells.Item(3, 2).Value = -26.49742398988303
$ws.Cells.Item(4, 2).Value = 222.4980796732398
$ws.Cells.Item(5, 2).Value = -12.68689882773273
$ws.Cells.Item(6, 2).Value = -5.148954746323284
$ws.Cells.Item(7, 2).Value = 4.803272209503419
$ws.Cells.Item(8, 2).Value = -11.95798969120987
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 347.4304745918557
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -37.51651331570361
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01789217"

# --- Sheet 15 ---
$ws = $wb.Worksheets.Item(15)
$ws.Cells.Item(2, 2).Value = 2.302661435332894
$ws.Cells.Item(3, 2).Value = -59.53963391068977
$ws.Cells.Item(4, 2).Value = 126.1348022820451
$ws.Cells.Item(5, 2).Value = 15.72443506663215
$ws.Cells.Item(6, 2).Value = -41.82043092192259
$ws.Cells.Item(7, 2).Value = -34.41132609085739
$ws.Cells.Item(8, 2).Value = -16.3172822240324
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 143.5846554722506
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 131.1620471247755
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01876225"

# --- Sheet 16 ---
$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(2, 2).Value = 66.91154225261148
$ws.Cells.Item(3, 2).Value = -547.5540553038156
$ws.Cells.Item(4, 2).Value = 512.726579926873
$ws.Cells.Item(5, 2).Value = 5.267730383574218
$ws.Cells.Item(6, 2).Value = 48.68316287234109
$ws.Cells.Item(7, 2).Value = -9.727524035362023
$ws.Cells.Item(8, 2).Value = -34.73428042418426
$ws.Cells.Item(9, 1).Value = "Commute_Trip"
$ws.Cells.Item(9, 2).Value = 158.899797710933
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = 56.80387636877367
$ws.Rows.Item(11).Delete()
$ws.Name = "summ01966017"
